$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Update existing row 92: observation_date -> 2025-09-30, DFII10 -> 1.8
$ws.Range("A92").Value = 45930
$ws.Range("B92").Value = 1.8

# Append new row 93: observation_date -> 2025-10-07, DFII10 -> 1.79
# (copy A92's date format down to A93 first, so the new cell keeps the
# same "yyyy-mm-dd" display style as the rest of column A)
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A93").Value = 45937
$ws.Range("B93").Value = 1.79

# Update the view: scroll so row 86 (col A) is the top-left cell, and the
# active selection moves to B94 (the first empty cell below the new data).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 86
$win.ScrollColumn = 1
$ws.Range("B94").Select()
